$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column I (9th column) - shifts old "Actual" data to column J
$ws.Range("I1").EntireColumn.Insert()

# New column I: Distribution channel code
$ws.Range("I1").Value = "Distribution channel code"
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"

# Match the header style (bold) used by row 1
$ws.Range("I1").Font.Bold = $true

# Set new column width (target OOXML stored width ~22.5546875 characters;
# this engine snaps ColumnWidth to 1/6-character steps, so 21.6667 is the
# closest achievable input and lands on the nearest representable width)
$ws.Range("I1").ColumnWidth = 21.6666666666667
